$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.642.06"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "2.038.80"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'231.68"
$ws.Range("E5").Value = "  -8.74%  "
$ws.Range("D6").Value = "'0.601"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "'55.37"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'57.06"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").Value = "'0.0752"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "2.338.36"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "'14.30"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("E15").Value = "  -5.75%  "
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "'5.15"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "2.027.10"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "36.756.71"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("E20").Value = "  +17.11%  "
$ws.Range("D21").Value = "'67.69"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "0.0₃0797"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("D23").Value = "'220.58"
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  -5.18%  "
$ws.Range("D27").Value = "'162.70"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "'8.73"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'18.89"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "'0.125"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("E31").Value = "  +4.94%  "
$ws.Range("D32").Value = "'0.117"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("D36").Value = "'4.28"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("D39").Value = "'5.79"
$ws.Range("E39").Value = "  +9.34%  "
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("D42").Value = "'4.39"
$ws.Range("E42").Value = "  +31.26%  "
$ws.Range("D43").Value = "1.474.33"
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("D44").Value = "'0.0942"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("D45").Value = "'93.90"
$ws.Range("E45").Value = "  +7.00%  "
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").Value = "'15.57"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "'6.95"
$ws.Range("E51").Value = "  +3.47%  "
